$p = $ppt.ActivePresentation
$newDate = "3/18/18"

function Set-DatePlaceholderText($shapes, $newText) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $shp = $shapes.Item($i)
        $isDate = $false
        try {
            if ($shp.PlaceholderFormat.Type -eq 16) {
                $isDate = $true
            }
        } catch {
            $isDate = $false
        }
        if ($isDate -and $shp.HasTextFrame) {
            $shp.TextFrame.TextRange.Text = $newText
        }
    }
}

# Slide master
$master = $p.SlideMaster
Set-DatePlaceholderText $master.Shapes $newDate

# Every slide layout (CustomLayout) under the master
for ($li = 1; $li -le $master.CustomLayouts.Count; $li++) {
    $layout = $master.CustomLayouts.Item($li)
    Set-DatePlaceholderText $layout.Shapes $newDate
}

# Notes master
$notesMaster = $p.NotesMaster
Set-DatePlaceholderText $notesMaster.Shapes $newDate
